$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2, 3, 4, 6, 7, 8, 9, 10, 11, 12, 15)

$data = @{
    2 = @(12.20844944890229, 9.295004803680674, 9.583523297650462, 32.70417578904911, 33.71147391099176, 15.69917436690564, 25.01088591646907, 10.54841625629543, 9.089983209259596, 11.52948286879961, 24.52109886189538)
    3 = @(11.92043621944675, 9.288058108571086, 9.561115437369303, 32.80792515795469, 33.85743155213324, 15.75150029690461, 25.11590079308582, 10.57190811446114, 8.873941003152284, 11.52016808732463, 24.61626759756905)
    4 = @(11.7411636790681, 9.284207861875208, 9.548688473687202, 32.87817757483656, 33.95565526138537, 15.78575269119057, 25.18447320768893, 10.58725009539191, 8.739167273784167, 11.51586019676357, 24.67905174443661)
    5 = @(11.66759631219151, 9.282744144049859, 9.543962944943745, 32.90845125411573, 33.99784056430345, 15.80024554252044, 25.21344708801874, 10.59373338961486, 8.683783719848025, 11.51446162233896, 24.70573054919035)
    6 = @(11.65535254397076, 9.282507480591244, 9.54319882676608, 32.91357749054478, 34.00497559795198, 15.8026843834823, 25.21832041865284, 10.59482392120195, 8.674561572035069, 11.51425100351185, 24.71022660112641)
    7 = @(11.74017346417642, 9.284187694147315, 9.548623367989531, 32.87857919713007, 33.95621545540672, 15.78594598096394, 25.18485978768379, 10.58733659411822, 8.738422125833349, 11.51583988736618, 24.67940711570959)
    8 = @(12.10970505286467, 9.292523994356504, 9.575523013354333, 32.73858811101527, 33.76001084956725, 15.71677586535016, 25.0462461030695, 10.55632604454393, 9.015974359460291, 11.52597943880442, 24.55301009069703)
    9 = @(12.81095996920235, 9.312124597656336, 9.638664706780823, 32.51611195702922, 33.44377251725862, 15.59795932205996, 24.80685981119372, 10.50277475229312, 9.540448200887148, 11.55696890234882, 24.33966789698332)
    10 = @(13.30678716051061, 9.328455008982937, 9.691143110101448, 32.38448487782367, 33.25352432190908, 15.52088391576788, 24.65069955472924, 10.46782525235185, 9.910022766831622, 11.58637531909217, 24.20397794468617)
    11 = @(13.52717194911321, 9.336291515383198, 9.716281524543485, 32.33153344792499, 33.17618157509717, 15.48803135875094, 24.58392822641826, 10.45287343622402, 10.07403549491323, 11.60116344850223, 24.14682233808906)
    12 = @(13.60980953554497, 9.339316578338043, 9.725977573686279, 32.31247943769814, 33.14822241573818, 15.47590810962366, 24.55925651235197, 10.44734723217255, 10.13549987500149, 11.60696332135856, 24.12583651262623)
    13 = @(13.59204952013892, 9.338662534871933, 9.723881585359576, 32.31653867074177, 33.15418473224548, 15.47850496260406, 24.56454274409842, 10.44853136962935, 10.12229184511937, 11.6057053719376, 24.13032692274425)
    14 = @(13.53398737687533, 9.336539242049692, 9.71707571474724, 32.32994586245657, 33.1738546944327, 15.48702761487448, 24.58188618367915, 10.45241607446103, 10.07910539435765, 11.60163661819536, 24.1450826359257)
    15 = @(13.49831409915424, 9.335246125906362, 9.712929769261683, 32.33828810755453, 33.18607633432497, 15.49228929545133, 24.59258937831467, 10.45481323075166, 10.05256712424418, 11.59917033503189, 24.15420661261366)
    16 = @(13.29227401413131, 9.327950987223963, 9.689525274167085, 32.38808489032686, 33.25876450589917, 15.52307533149038, 24.6551490396986, 10.4688214050226, 9.899216856254601, 11.5854370224242, 24.2078052126273)
    17 = @(13.16449655419792, 9.323579319577492, 9.675487830019444, 32.42040876493578, 33.30571727426432, 15.54252715826989, 24.69461987622004, 10.47765717422828, 9.80405013964643, 11.5773714046553, 24.24185719506391)
    18 = @(13.09051980952198, 9.321103294504146, 9.667533332312084, 32.43965246495414, 33.33358908688117, 15.55392332053703, 24.71772408703898, 10.482828426264, 9.748928926370024, 11.57286527218047, 24.26187315240779)
    19 = @(13.06539203085998, 9.320271590132306, 9.664860753995079, 32.44627995954042, 33.34317453565477, 15.55781759985206, 24.72561576941282, 10.48459464855793, 9.730201588113594, 11.57136250371203, 24.26872407542159)
    20 = @(13.17814917750748, 9.324040720179408, 9.67696981330239, 32.41690035633241, 33.30062942136799, 15.54043495536065, 24.69037657426222, 10.47670736787005, 9.814220921170271, 11.57821625974799, 24.23818777982483)
    21 = @(13.55106437891708, 9.337161351500114, 9.719070012478538, 32.32598076063773, 33.16804104107975, 15.48451569706678, 24.57677535776924, 10.45127136242806, 10.0918081711853, 11.60282630818599, 24.14073066943735)
    22 = @(13.78999168499723, 9.34607136112426, 9.747612119516168, 32.2723746157477, 33.08913495147745, 15.44981852165858, 24.50610407297684, 10.43543836765565, 10.26945410433152, 11.62007425649992, 24.08087101531405)
    23 = @(13.66293329554516, 9.341285638609547, 9.732286486029629, 32.30045267826468, 33.13053786275466, 15.46816797957595, 24.54349580315371, 10.44381651020955, 10.17500275410876, 11.61076322813021, 24.1124682280822)
    24 = @(13.17197842737116, 9.323832004688267, 9.676299447514811, 32.41848445128886, 33.30292690506726, 15.54138017663412, 24.69229368800359, 10.47713649043317, 9.809623979482824, 11.57783389280813, 24.23984535498283)
    25 = @(12.62430719958817, 9.306479273019026, 9.620495235058552, 32.57071517875386, 33.52195426877424, 15.62830504112922, 24.86815374097317, 10.50277475229312, 9.401080866068156, 11.54740949344822, 24.39368613451557)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $rowVals[$i]
    }
}